# Automatic update of files.
# - Round the Ost/Nord (easting/northing) coordinates in row 2 to whole numbers
# - Drop the Starttid/Sluttid ("00:00") time values, leaving those cells empty

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 506128
$ws.Range("R2").Value = 6932995

$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
